$d = $word.ActiveDocument

# Change 1: merge " " + "Lahav" into " Lahav" (removes proofErr spell check wrapper)
$d.Content.Find.Execute(" Lahav", $true, $false, $false, $false, $false, $true, 1, $false, " Lahav", 2) | Out-Null
